$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change A92's number format to match the "in-progress" rows (YYYY-MM-DD HH:MM:SS)
$ws.Range("A92").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add new row 93 with the next day's data
$ws.Range("A93").Value = 45832
$ws.Range("A93").NumberFormat = "YYYY-MM-DD"
$ws.Range("B93").Value = 394
$ws.Range("C93").Value = 395
$ws.Range("D93").Value = 400
